$wb = $excel.ActiveWorkbook

# --- Backlog_Produto: update task-estimate hours (the prototipo addition
#     increases several "Construção"/"Transição" task estimates) ---
$ws = $wb.Worksheets.Item("Backlog_Produto")

$ws.Range("D15").Value = 15
$ws.Range("D17").Value = 6
$ws.Range("D18").Value = 6
$ws.Range("D19").Value = 4
$ws.Range("D20").Value = 4

$ws.Range("D23").Value = 15
$ws.Range("D24").Value = 12
$ws.Range("D25").Value = 15
$ws.Range("D26").Value = 12
$ws.Range("D28").Value = 10

$ws.Range("D33").Value = 4
$ws.Range("D34").Value = 4

$ws.Range("D37").Value = 6
$ws.Range("D39").Value = 10
$ws.Range("D40").Value = 4
$ws.Range("D41").Value = 4
$ws.Range("D43").Value = 4
$ws.Range("D44").Value = 4

# --- Restore/adjust the active selections on the sheets that had their
#     cursor position changed by the author ---
$wsPlan = $wb.Worksheets.Item("Planejamento")
$wsPlan.Range("A9").Select()

$ws.Activate()
$ws.Range("F16").Select()

$excel.Calculate()
